# DOMA-3382 add payment status to excel export
# Adds a "status" column (between "transaction" and "amount") to the
# payments export template, and adds a block of empty, bordered rows
# below the 3-row header/sample area (used by the export engine as a
# repeating/formatting area).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert a new column at G, shifting the old "amount" column (G) to H.
#    The newly freed column G becomes the "status" column.
# ---------------------------------------------------------------------
$ws.Columns.Item(7).Insert()

$ws.Range("G1").Value = "{d.i18n.status}"
$ws.Range("G2").Value = "{d.objs[I].status}"
$ws.Range("G3").Value = "{d.objs[i+1].status}"

# Give the new header/body cells the same look & feel as their neighbours.
$ws.Range("G1").Font.Size = $ws.Range("F1").Font.Size
$ws.Range("G1").Font.Name = $ws.Range("F1").Font.Name
$ws.Range("G1").Font.Bold = $ws.Range("F1").Font.Bold
$ws.Range("G1").Interior.Color = $ws.Range("F1").Interior.Color
$ws.Range("G2:G3").Font.Size = $ws.Range("F2").Font.Size
$ws.Range("G2:G3").Font.Name = $ws.Range("F2").Font.Name
$ws.Range("G2:G3").Interior.Color = $ws.Range("F2").Interior.Color

# ---------------------------------------------------------------------
# 2. Append 7 new, empty, formatted rows (4-10) below the header block.
#    Row 4 forms the top edge of a bordered box, rows 5-9 are the body
#    (only left/right edges), and row 10 is the bottom edge. The box is
#    drawn in a light grey, except the very top edge which is black.
# ---------------------------------------------------------------------
$blackColor = 0        # RGB(0,0,0)
$grayColor  = 11184810 # RGB(170,170,170) = #AAAAAA
$whiteColor = 16777215 # RGB(255,255,255)

$ws.Rows.Item(4).Resize(7).Insert()

for ($r = 4; $r -le 10; $r++) {
    $ws.Range("A" + $r + ":H" + $r).Interior.Color = $whiteColor
    $ws.Range("A" + $r + ":H" + $r).RowHeight = 13.65
}

# Row 4 - top edge of the box (black top border, grey left/right edges)
$ws.Range("A4:H4").Borders.Item(8).LineStyle = 1
$ws.Range("A4:H4").Borders.Item(8).Color = $blackColor
$ws.Range("A4").Borders.Item(7).LineStyle = 1
$ws.Range("A4").Borders.Item(7).Color = $grayColor
$ws.Range("H4").Borders.Item(10).LineStyle = 1
$ws.Range("H4").Borders.Item(10).Color = $grayColor

# Rows 5-9 - body, only the outer left/right edges are drawn
$ws.Range("A5:A9").Borders.Item(7).LineStyle = 1
$ws.Range("A5:A9").Borders.Item(7).Color = $grayColor
$ws.Range("H5:H9").Borders.Item(10).LineStyle = 1
$ws.Range("H5:H9").Borders.Item(10).Color = $grayColor

# Row 10 - bottom edge of the box (grey bottom/left/right edges)
$ws.Range("A10:H10").Borders.Item(9).LineStyle = 1
$ws.Range("A10:H10").Borders.Item(9).Color = $grayColor
$ws.Range("A10").Borders.Item(7).LineStyle = 1
$ws.Range("A10").Borders.Item(7).Color = $grayColor
$ws.Range("H10").Borders.Item(10).LineStyle = 1
$ws.Range("H10").Borders.Item(10).Color = $grayColor
